$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column K (Description), shifting it to L
$ws.Range("K1").EntireColumn.Insert()

# New column K header and values
$ws.Range("K1").Value = "Graphic File"
$ws.Range("K2:K7").Value = "None"
$ws.Range("K2").Value = "boeing_787-8_cropped"
$ws.Range("K5").Value = "B737_100"

# Set selection to K6
$ws.Range("K6").Select()
